$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 54302.855
$ws.Range("I40").Value = 10450.909
$ws.Range("J40").Value = 102540
$ws.Range("K40").Value = 10450.909
$ws.Range("L40").Value = 102540
$ws.Range("M40").Value = -10275.909
$ws.Range("N40").Value = -102890

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1624
$ws.Range("I45").Value = 1286.2
$ws.Range("J45").Value = 1961.8
$ws.Range("K45").Value = 1286.2
$ws.Range("L45").Value = 1961.8
$ws.Range("M45").Value = -909.2
$ws.Range("N45").Value = -2715.8

$ws.Range("H74").Value = 1391.6666
$ws.Range("I74").Value = 842.8570999999999
$ws.Range("J74").Value = 2160
$ws.Range("K74").Value = 842.8570999999999
$ws.Range("L74").Value = 2160
$ws.Range("M74").Value = 31.14290000000005
$ws.Range("N74").Value = -3908

$ws.Range("H77").Value = 1391.6666
$ws.Range("I77").Value = 842.8570999999999
$ws.Range("J77").Value = 2160
$ws.Range("K77").Value = 4214.2855
$ws.Range("L77").Value = 10800
$ws.Range("M77").Value = 153.7145
$ws.Range("N77").Value = -19536

$ws.Range("H110").Value = 2250.0833
$ws.Range("I110").Value = 1400.1
$ws.Range("J110").Value = 6500
$ws.Range("K110").Value = 1400.1
$ws.Range("L110").Value = 6500
$ws.Range("M110").Value = 644.9000000000001
$ws.Range("N110").Value = -10590

$ws.Range("H122").Value = 1589.3334
$ws.Range("I122").Value = 1593.1052
$ws.Range("J122").Value = 1553.5
$ws.Range("K122").Value = 4779.3156
$ws.Range("L122").Value = 4660.5
$ws.Range("M122").Value = -2329.3156
$ws.Range("N122").Value = -9560.5

$ws.Range("H132").Value = 2746.25
$ws.Range("I132").Value = 2045.6
$ws.Range("J132").Value = 6249.5
$ws.Range("K132").Value = 6136.799999999999
$ws.Range("L132").Value = 18748.5
$ws.Range("M132").Value = -3606.799999999999
$ws.Range("N132").Value = -23808.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1738.6923
$ws.Range("I86").Value = 1812.6666
$ws.Range("J86").Value = 1428
$ws.Range("K86").Value = 1812.6666
$ws.Range("L86").Value = 1428
$ws.Range("M86").Value = -689.6666
$ws.Range("N86").Value = -3674

$ws.Range("H89").Value = 1738.6923
$ws.Range("I89").Value = 1812.6666
$ws.Range("J89").Value = 1428
$ws.Range("K89").Value = 9063.333000000001
$ws.Range("L89").Value = 7140
$ws.Range("M89").Value = -3447.333000000001
$ws.Range("N89").Value = -18372

$ws.Range("H107").Value = 1700
$ws.Range("I107").Value = 1266.6666
$ws.Range("J107").Value = 3000
$ws.Range("K107").Value = 1266.6666
$ws.Range("L107").Value = 3000
$ws.Range("M107").Value = 653.3334
$ws.Range("N107").Value = -6840

$ws.Range("H134").Value = 1995.9062
$ws.Range("I134").Value = 1862.7727
$ws.Range("J134").Value = 2288.8
$ws.Range("K134").Value = 5588.3181
$ws.Range("L134").Value = 6866.400000000001
$ws.Range("M134").Value = -3053.3181
$ws.Range("N134").Value = -11936.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 2912
$ws.Range("I11").Value = 1000
$ws.Range("J11").Value = 3294.4
$ws.Range("K11").Value = 1000
$ws.Range("L11").Value = 3294.4
$ws.Range("M11").Value = -860
$ws.Range("N11").Value = -3574.4

$ws.Range("H16").Value = 3945.3076
$ws.Range("I16").Value = 3396
$ws.Range("J16").Value = 4288.625
$ws.Range("K16").Value = 3396
$ws.Range("L16").Value = 4288.625
$ws.Range("M16").Value = -3109
$ws.Range("N16").Value = -4862.625

$ws.Range("H31").Value = 1442.1296
$ws.Range("I31").Value = 841.0769
$ws.Range("J31").Value = 2000.25
$ws.Range("K31").Value = 841.0769
$ws.Range("L31").Value = 2000.25
$ws.Range("M31").Value = -546.0769
$ws.Range("N31").Value = -2590.25

$ws.Range("H34").Value = 1442.1296
$ws.Range("I34").Value = 841.0769
$ws.Range("J34").Value = 2000.25
$ws.Range("K34").Value = 841.0769
$ws.Range("L34").Value = 2000.25
$ws.Range("M34").Value = -639.0769
$ws.Range("N34").Value = -2404.25

$ws.Range("H113").Value = 3945.3076
$ws.Range("I113").Value = 3396
$ws.Range("J113").Value = 4288.625
$ws.Range("K113").Value = 3396
$ws.Range("L113").Value = 4288.625
$ws.Range("M113").Value = -1226
$ws.Range("N113").Value = -8628.625

$ws.Range("H132").Value = 2162.139
$ws.Range("I132").Value = 1590.5172
$ws.Range("J132").Value = 4530.2856
$ws.Range("K132").Value = 4771.5516
$ws.Range("L132").Value = 13590.8568
$ws.Range("M132").Value = -2241.5516
$ws.Range("N132").Value = -18650.8568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 4383.3335
$ws.Range("J51").Value = 4806.25
$ws.Range("L51").Value = 14418.75
$ws.Range("N51").Value = -15338.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1137.0741
$ws.Range("I107").Value = 1610.4375
$ws.Range("J107").Value = 448.54544
$ws.Range("K107").Value = 1610.4375
$ws.Range("L107").Value = 448.54544
$ws.Range("M107").Value = 309.5625
$ws.Range("N107").Value = -4288.54544

$ws.Range("H113").Value = 8260.733
$ws.Range("I113").Value = 9759.166999999999
$ws.Range("J113").Value = 2267
$ws.Range("K113").Value = 9759.166999999999
$ws.Range("L113").Value = 2267
$ws.Range("M113").Value = -7589.166999999999
$ws.Range("N113").Value = -6607

$ws.Range("H122").Value = 2405.1562
$ws.Range("I122").Value = 2090.6538
$ws.Range("J122").Value = 3768
$ws.Range("K122").Value = 6271.9614
$ws.Range("L122").Value = 11304
$ws.Range("M122").Value = -3821.9614
$ws.Range("N122").Value = -16204

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3200.8
$ws.Range("I61").Value = 3200.8
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3200.8
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2998.8
$ws.Range("N61").ClearContents()

$ws.Range("H68").Value = 13925
$ws.Range("I68").Value = 17933.334
$ws.Range("J68").Value = 1900
$ws.Range("K68").Value = 17933.334
$ws.Range("L68").Value = 1900
$ws.Range("M68").Value = -17184.334
$ws.Range("N68").Value = -3398

$ws.Range("H71").Value = 13925
$ws.Range("I71").Value = 17933.334
$ws.Range("J71").Value = 1900
$ws.Range("K71").Value = 89666.67
$ws.Range("L71").Value = 9500
$ws.Range("M71").Value = -85922.67
$ws.Range("N71").Value = -16988

$ws.Range("H112").Value = 18387
$ws.Range("J112").Value = 18387
$ws.Range("L112").Value = 18387
$ws.Range("N112").Value = -21341

$ws.Range("H113").Value = 3200.8
$ws.Range("I113").Value = 3200.8
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3200.8
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1030.8
$ws.Range("N113").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 39375
$ws.Range("J103").Value = 39375
$ws.Range("L103").Value = 39375
$ws.Range("N103").Value = -41719

$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180

$ws.Range("H114").Value = 40000
$ws.Range("J114").Value = 40000
$ws.Range("L114").Value = 40000
$ws.Range("N114").Value = -48678

$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
